$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H2").Value = 1092.625
$ws_ALC.Range("I2").Value = 391.57144
$ws_ALC.Range("J2").Value = 6000
$ws_ALC.Range("K2").Value = 391.57144
$ws_ALC.Range("L2").Value = 6000
$ws_ALC.Range("M2").Value = -278.57144
$ws_ALC.Range("N2").Value = -6226

$ws_ALC.Range("H15").Value = 326.42856
$ws_ALC.Range("I15").Value = 326.42856
$ws_ALC.Range("K15").Value = 979.28568
$ws_ALC.Range("M15").Value = -810.28568

$ws_ALC.Range("H33").Value = 217.45454
$ws_ALC.Range("I33").Value = 169.22223
$ws_ALC.Range("K33").Value = 169.22223
$ws_ALC.Range("M33").Value = 59.77777

$ws_ALC.Range("H100").Value = 2272.2856
$ws_ALC.Range("I100").Value = 2272.2856
$ws_ALC.Range("K100").Value = 2272.2856
$ws_ALC.Range("M100").Value = -1731.2856

$ws_ALC.Range("H103").Value = 2416.7368
$ws_ALC.Range("I103").Value = 580
$ws_ALC.Range("J103").Value = 4457.5557
$ws_ALC.Range("K103").Value = 1740
$ws_ALC.Range("L103").Value = 13372.6671
$ws_ALC.Range("M103").Value = -1154
$ws_ALC.Range("N103").Value = -14544.6671

$ws_ALC.Range("H113").Value = 9037.166999999999
$ws_ALC.Range("I113").Value = 8844.1
$ws_ALC.Range("K113").Value = 8844.1
$ws_ALC.Range("M113").Value = -5590.1

$ws_ALC.Range("H137").Value = 1546
$ws_ALC.Range("I137").Value = 1546
$ws_ALC.Range("K137").Value = 4638
$ws_ALC.Range("M137").Value = -2088

$ws_ARM.Range("H45").Value = 3099.5715
$ws_ARM.Range("I45").Value = 1378.7142
$ws_ARM.Range("K45").Value = 1378.7142
$ws_ARM.Range("M45").Value = -1001.7142

$ws_ARM.Range("H56").Value = 35000
$ws_ARM.Range("I56").Value = 0
$ws_ARM.Range("J56").Value = 35000
$ws_ARM.Range("K56").Value = 0
$ws_ARM.Range("L56").Value = 35000
$ws_ARM.Range("M56").ClearContents()
$ws_ARM.Range("N56").Value = -36484

$ws_ARM.Range("H61").Value = 591
$ws_ARM.Range("I61").Value = 372
$ws_ARM.Range("K61").Value = 372
$ws_ARM.Range("M61").Value = -160

$ws_ARM.Range("H102").Value = 3144.6
$ws_ARM.Range("I102").Value = 3144.6
$ws_ARM.Range("J102").Value = 0
$ws_ARM.Range("K102").Value = 3144.6
$ws_ARM.Range("L102").Value = 0
$ws_ARM.Range("M102").Value = -1522.6
$ws_ARM.Range("N102").ClearContents()

$ws_ARM.Range("H132").Value = 1564.5
$ws_ARM.Range("I132").Value = 1550
$ws_ARM.Range("J132").Value = 1695
$ws_ARM.Range("K132").Value = 4650
$ws_ARM.Range("L132").Value = 5085
$ws_ARM.Range("M132").Value = -2120
$ws_ARM.Range("N132").Value = -10145

$ws_ARM.Range("H136").Value = 591
$ws_ARM.Range("I136").Value = 372
$ws_ARM.Range("K136").Value = 1116
$ws_ARM.Range("M136").Value = 1434

$ws_BSM.Range("H105").Value = 2626.8125
$ws_BSM.Range("I105").Value = 2626.8125
$ws_BSM.Range("K105").Value = 2626.8125
$ws_BSM.Range("M105").Value = -879.8125

$ws_CRP.Range("H10").Value = 2122.6
$ws_CRP.Range("I10").Value = 153.25
$ws_CRP.Range("J10").Value = 10000
$ws_CRP.Range("K10").Value = 153.25
$ws_CRP.Range("L10").Value = 10000
$ws_CRP.Range("M10").Value = -14.25
$ws_CRP.Range("N10").Value = -10278

$ws_CRP.Range("H31").Value = 3624.5557
$ws_CRP.Range("I31").Value = 2939.1667
$ws_CRP.Range("J31").Value = 4995.3335
$ws_CRP.Range("K31").Value = 2939.1667
$ws_CRP.Range("L31").Value = 4995.3335
$ws_CRP.Range("M31").Value = -2644.1667
$ws_CRP.Range("N31").Value = -5585.3335

$ws_CRP.Range("H34").Value = 3624.5557
$ws_CRP.Range("I34").Value = 2939.1667
$ws_CRP.Range("J34").Value = 4995.3335
$ws_CRP.Range("K34").Value = 2939.1667
$ws_CRP.Range("L34").Value = 4995.3335
$ws_CRP.Range("M34").Value = -2737.1667
$ws_CRP.Range("N34").Value = -5399.3335

$ws_CRP.Range("H38").Value = 6000
$ws_CRP.Range("I38").Value = 0
$ws_CRP.Range("K38").Value = 0
$ws_CRP.Range("M38").ClearContents()

$ws_CRP.Range("H46").Value = 6000
$ws_CRP.Range("I46").Value = 0
$ws_CRP.Range("K46").Value = 0
$ws_CRP.Range("M46").ClearContents()

$ws_CRP.Range("H58").Value = 0
$ws_CRP.Range("I58").Value = 0
$ws_CRP.Range("J58").Value = 0
$ws_CRP.Range("K58").Value = 0
$ws_CRP.Range("L58").Value = 0
$ws_CRP.Range("M58").ClearContents()
$ws_CRP.Range("N58").ClearContents()

$ws_CRP.Range("H108").Value = 41661.668
$ws_CRP.Range("I108").Value = 0
$ws_CRP.Range("J108").Value = 41661.668
$ws_CRP.Range("K108").Value = 0
$ws_CRP.Range("L108").Value = 41661.668
$ws_CRP.Range("M108").ClearContents()
$ws_CRP.Range("N108").Value = -49341.668

$ws_CRP.Range("H122").Value = 2736.1667
$ws_CRP.Range("I122").Value = 2796.4
$ws_CRP.Range("K122").Value = 8389.200000000001
$ws_CRP.Range("M122").Value = -5939.200000000001

$ws_CRP.Range("H132").Value = 2501.4666
$ws_CRP.Range("I132").Value = 2501.4666
$ws_CRP.Range("K132").Value = 7504.399800000001
$ws_CRP.Range("M132").Value = -4974.399800000001

$ws_CRP.Range("H134").Value = 1501.6666
$ws_CRP.Range("I134").Value = 1501.6666
$ws_CRP.Range("J134").Value = 0
$ws_CRP.Range("K134").Value = 4504.9998
$ws_CRP.Range("L134").Value = 0
$ws_CRP.Range("M134").Value = -1969.9998
$ws_CRP.Range("N134").ClearContents()

$ws_CRP.Range("H136").Value = 0
$ws_CRP.Range("I136").Value = 0
$ws_CRP.Range("J136").Value = 0
$ws_CRP.Range("K136").Value = 0
$ws_CRP.Range("L136").Value = 0
$ws_CRP.Range("M136").ClearContents()
$ws_CRP.Range("N136").ClearContents()

$ws_CUL.Range("H81").Value = 850
$ws_CUL.Range("J81").Value = 1500
$ws_CUL.Range("L81").Value = 4500
$ws_CUL.Range("N81").Value = -6746

$ws_CUL.Range("H84").Value = 850
$ws_CUL.Range("J84").Value = 1500
$ws_CUL.Range("L84").Value = 13500
$ws_CUL.Range("N84").Value = -24732

$ws_CUL.Range("H86").Value = 1263.3334
$ws_CUL.Range("I86").Value = 1116
$ws_CUL.Range("K86").Value = 3348
$ws_CUL.Range("M86").Value = -2162

$ws_CUL.Range("H89").Value = 1263.3334
$ws_CUL.Range("I89").Value = 1116
$ws_CUL.Range("K89").Value = 10044
$ws_CUL.Range("M89").Value = -4116

$ws_CUL.Range("H114").Value = 459.375
$ws_CUL.Range("I114").Value = 484.66666
$ws_CUL.Range("J114").Value = 444.2
$ws_CUL.Range("K114").Value = 1453.99998
$ws_CUL.Range("L114").Value = 1332.6
$ws_CUL.Range("M114").Value = 1800.00002
$ws_CUL.Range("N114").Value = -7840.6

$ws_CUL.Range("H132").Value = 2252.375
$ws_CUL.Range("I132").Value = 2145.5715
$ws_CUL.Range("J132").Value = 3000
$ws_CUL.Range("K132").Value = 19310.1435
$ws_CUL.Range("L132").Value = 27000
$ws_CUL.Range("M132").Value = -16780.1435
$ws_CUL.Range("N132").Value = -32060

$ws_GSM.Range("H102").Value = 2472.8333
$ws_GSM.Range("I102").Value = 2360.5
$ws_GSM.Range("J102").Value = 2697.5
$ws_GSM.Range("K102").Value = 2360.5
$ws_GSM.Range("L102").Value = 2697.5
$ws_GSM.Range("M102").Value = -738.5
$ws_GSM.Range("N102").Value = -5941.5

$ws_GSM.Range("H132").Value = 2404.1428
$ws_GSM.Range("I132").Value = 2403.3845
$ws_GSM.Range("K132").Value = 7210.1535
$ws_GSM.Range("M132").Value = -4680.1535

$ws_GSM.Range("H134").Value = 15000
$ws_GSM.Range("J134").Value = 15000
$ws_GSM.Range("L134").Value = 45000
$ws_GSM.Range("N134").Value = -50070

$ws_LTW.Range("H9").Value = 295.8
$ws_LTW.Range("I9").Value = 342
$ws_LTW.Range("J9").Value = 111
$ws_LTW.Range("K9").Value = 342
$ws_LTW.Range("L9").Value = 111
$ws_LTW.Range("M9").Value = -118
$ws_LTW.Range("N9").Value = -559

$ws_LTW.Range("H22").Value = 2549.5
$ws_LTW.Range("I22").Value = 499
$ws_LTW.Range("K22").Value = 499
$ws_LTW.Range("M22").Value = -204

$ws_LTW.Range("H27").Value = 2549.5
$ws_LTW.Range("I27").Value = 499
$ws_LTW.Range("K27").Value = 499
$ws_LTW.Range("M27").Value = -392

$ws_LTW.Range("H33").Value = 69420
$ws_LTW.Range("I33").Value = 0
$ws_LTW.Range("K33").Value = 0
$ws_LTW.Range("M33").ClearContents()

$ws_LTW.Range("H35").Value = 1039.75
$ws_LTW.Range("I35").Value = 1039.75
$ws_LTW.Range("K35").Value = 1039.75
$ws_LTW.Range("M35").Value = -703.75

$ws_LTW.Range("H132").Value = 1844.2858
$ws_LTW.Range("I132").Value = 1844.2858
$ws_LTW.Range("J132").Value = 0
$ws_LTW.Range("K132").Value = 5532.857400000001
$ws_LTW.Range("L132").Value = 0
$ws_LTW.Range("M132").Value = -3002.857400000001
$ws_LTW.Range("N132").ClearContents()

$ws_LTW.Range("H135").Value = 65000
$ws_LTW.Range("J135").Value = 65000
$ws_LTW.Range("L135").Value = 65000
$ws_LTW.Range("N135").Value = -75140

$ws_WVR.Range("H122").Value = 0
$ws_WVR.Range("J122").Value = 0
$ws_WVR.Range("L122").Value = 0
$ws_WVR.Range("N122").ClearContents()

$ws_WVR.Range("H135").Value = 47497.5
$ws_WVR.Range("J135").Value = 47497.5
$ws_WVR.Range("L135").Value = 47497.5
$ws_WVR.Range("N135").Value = -57637.5
